$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 11111851
$ws.Range("I43").Value = 833.3333
$ws.Range("K43").Value = 833.3333
$ws.Range("M43").Value = -764.3333

$ws.Range("H86").Value = 5007.5625
$ws.Range("I86").Value = 2382.3333
$ws.Range("K86").Value = 2382.3333
$ws.Range("M86").Value = -1259.3333

$ws.Range("H89").Value = 5007.5625
$ws.Range("I89").Value = 2382.3333
$ws.Range("K89").Value = 11911.6665
$ws.Range("M89").Value = -6295.666499999999

$ws.Range("H100").Value = 2332.6667
$ws.Range("I100").Value = 1000
$ws.Range("K100").Value = 1000
$ws.Range("M100").Value = -459

$ws.Range("H126").Value = 29996
$ws.Range("J126").Value = 29996
$ws.Range("L126").Value = 29996
$ws.Range("N126").Value = -39876

$ws.Range("H138").Value = 2099.16
$ws.Range("I138").Value = 1554.4546
$ws.Range("J138").Value = 2166.4832
$ws.Range("K138").Value = 4663.3638
$ws.Range("L138").Value = 6499.4496
$ws.Range("M138").Value = 476.6361999999999
$ws.Range("N138").Value = -16779.4496

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1093.3529
$ws.Range("I61").Value = 868.8461
$ws.Range("J61").Value = 1823
$ws.Range("K61").Value = 868.8461
$ws.Range("L61").Value = 1823
$ws.Range("M61").Value = -656.8461
$ws.Range("N61").Value = -2247

$ws.Range("H74").Value = 811.0417
$ws.Range("I74").Value = 813.2105
$ws.Range("K74").Value = 813.2105
$ws.Range("M74").Value = 60.78949999999998

$ws.Range("H77").Value = 811.0417
$ws.Range("I77").Value = 813.2105
$ws.Range("K77").Value = 4066.0525
$ws.Range("M77").Value = 301.9474999999998

$ws.Range("H136").Value = 1093.3529
$ws.Range("I136").Value = 868.8461
$ws.Range("J136").Value = 1823
$ws.Range("K136").Value = 2606.5383
$ws.Range("L136").Value = 5469
$ws.Range("M136").Value = -56.53830000000016
$ws.Range("N136").Value = -10569

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1340.2222
$ws.Range("I20").Value = 1136.7858
$ws.Range("J20").Value = 2052.25
$ws.Range("K20").Value = 1136.7858
$ws.Range("L20").Value = 2052.25
$ws.Range("M20").Value = -889.7858000000001
$ws.Range("N20").Value = -2546.25

$ws.Range("H86").Value = 3758.182
$ws.Range("I86").Value = 3888.8823
$ws.Range("J86").Value = 3313.8
$ws.Range("K86").Value = 3888.8823
$ws.Range("L86").Value = 3313.8
$ws.Range("M86").Value = -2765.8823
$ws.Range("N86").Value = -5559.8

$ws.Range("H89").Value = 3758.182
$ws.Range("I89").Value = 3888.8823
$ws.Range("J89").Value = 3313.8
$ws.Range("K89").Value = 19444.4115
$ws.Range("L89").Value = 16569
$ws.Range("M89").Value = -13828.4115
$ws.Range("N89").Value = -27801

$ws.Range("H132").Value = 5007500
$ws.Range("J132").Value = 5007500
$ws.Range("L132").Value = 5007500
$ws.Range("N132").Value = -5017620

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 738
$ws.Range("I10").Value = 738
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 738
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -599
$ws.Range("N10").ClearContents()

$ws.Range("H62").Value = 4764395.5
$ws.Range("I62").Value = 2551.3416
$ws.Range("K62").Value = 2551.3416
$ws.Range("M62").Value = -1927.3416

$ws.Range("H65").Value = 4764395.5
$ws.Range("I65").Value = 2551.3416
$ws.Range("K65").Value = 12756.708
$ws.Range("M65").Value = -9636.708000000001

$ws.Range("H98").Value = 31480
$ws.Range("J98").Value = 31480
$ws.Range("L98").Value = 31480
$ws.Range("N98").Value = -35972

$ws.Range("H122").Value = 950.3333
$ws.Range("J122").Value = 938
$ws.Range("L122").Value = 2814
$ws.Range("N122").Value = -7714

$ws.Range("H127").Value = 0
$ws.Range("I127").Value = 0
$ws.Range("K127").Value = 0
$ws.Range("M127").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 657350.2
$ws.Range("I4").Value = 733487.1
$ws.Range("J4").Value = 605438.6
$ws.Range("K4").Value = 2200461.3
$ws.Range("L4").Value = 1816315.8
$ws.Range("M4").Value = -2200349.3
$ws.Range("N4").Value = -1816539.8

$ws.Range("H70").Value = 4558.8237
$ws.Range("I70").Value = 1725
$ws.Range("J70").Value = 5430.769
$ws.Range("K70").Value = 5175
$ws.Range("L70").Value = 16292.307
$ws.Range("M70").Value = -4860
$ws.Range("N70").Value = -16922.307

$ws.Range("H73").Value = 4558.8237
$ws.Range("I73").Value = 1725
$ws.Range("J73").Value = 5430.769
$ws.Range("K73").Value = 5175
$ws.Range("L73").Value = 16292.307
$ws.Range("M73").Value = -4083
$ws.Range("N73").Value = -18476.307

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 3000
$ws.Range("J40").Value = 3000
$ws.Range("L40").Value = 3000
$ws.Range("N40").Value = -3302

$ws.Range("H52").Value = 19975
$ws.Range("J52").Value = 19975
$ws.Range("L52").Value = 19975
$ws.Range("N52").Value = -20493

$ws.Range("H102").Value = 1782.9
$ws.Range("J102").Value = 1673.4286
$ws.Range("L102").Value = 1673.4286
$ws.Range("N102").Value = -4917.4286

$ws.Range("H122").Value = 2886.762
$ws.Range("I122").Value = 1762.2
$ws.Range("K122").Value = 5286.6
$ws.Range("M122").Value = -2836.6

$ws.Range("H132").Value = 2312.125
$ws.Range("I132").Value = 1888.2727
$ws.Range("J132").Value = 3244.6
$ws.Range("K132").Value = 5664.8181
$ws.Range("L132").Value = 9733.799999999999
$ws.Range("M132").Value = -3134.8181
$ws.Range("N132").Value = -14793.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 637.9
$ws.Range("I22").Value = 317.9
$ws.Range("J22").Value = 957.9
$ws.Range("K22").Value = 317.9
$ws.Range("L22").Value = 957.9
$ws.Range("M22").Value = -22.89999999999998
$ws.Range("N22").Value = -1547.9

$ws.Range("H27").Value = 637.9
$ws.Range("I27").Value = 317.9
$ws.Range("J27").Value = 957.9
$ws.Range("K27").Value = 317.9
$ws.Range("L27").Value = 957.9
$ws.Range("M27").Value = -210.9
$ws.Range("N27").Value = -1171.9

$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()

$ws.Range("H122").Value = 31253686
$ws.Range("I122").Value = 62503130
$ws.Range("J122").Value = 4245.75
$ws.Range("K122").Value = 187509390
$ws.Range("L122").Value = 12737.25
$ws.Range("M122").Value = -187506940
$ws.Range("N122").Value = -17637.25

$ws.Range("H132").Value = 69320.07000000001
$ws.Range("I132").Value = 2233.111
$ws.Range("J132").Value = 169950.5
$ws.Range("K132").Value = 6699.333
$ws.Range("L132").Value = 509851.5
$ws.Range("M132").Value = -4169.333
$ws.Range("N132").Value = -514911.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 440.4889
$ws.Range("I136").Value = 376.2963
$ws.Range("J136").Value = 536.7778
$ws.Range("K136").Value = 1128.8889
$ws.Range("L136").Value = 1610.3334
$ws.Range("M136").Value = 1421.1111
$ws.Range("N136").Value = -6710.3334
